$d = $word.ActiveDocument

$replacements = @(
    @("23×38=874", "54×38=2052"),
    @("51×40=2040", "88×40=3520"),
    @("91×11=1001", "55×69=3795"),
    @("18×56=1008", "21×45=945"),
    @("51×42=2142", "15×85=1275"),
    @("32×88=2816", "56×44=2464"),
    @("51×84=4284", "38×90=3420"),
    @("95×56=5320", "97×72=6984"),
    @("43×57=2451", "74×35=2590"),
    @("77×15=1155", "20×65=1300"),
    @("17×20=340", "89×37=3293"),
    @("91×62=5642", "86×91=7826"),
    @("56×23=1288", "53×44=2332"),
    @("38×15=570", "64×13=832"),
    @("59×49=2891", "64×12=768"),
    @("46×54=2484", "29×88=2552"),
    @("57×85=4845", "87×38=3306"),
    @("80×35=2800", "34×20=680"),
    @("94×81=7614", "35×64=2240"),
    @("77×29=2233", "41×13=533"),
    @("41×19=779", "57×11=627"),
    @("43×49=2107", "11×98=1078"),
    @("48×12=576", "66×25=1650"),
    @("40×75=3000", "62×42=2604"),
    @("82×44=3608", "59×43=2537")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
